$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-12 (columns A-H)
$data = @(
    @("Poland", 37950802, 312679, 121.37, "Europe", "Central Europe", "Polish", "PLN (Polish złoty)"),
    @("Finland", 5530719, 338424, 16.34, "Europe", "Northern Europe", "Finnish, Swedish", "EUR (Euro)"),
    @("Iceland", 366425, 103000, 3.56, "Europe", "Northern Europe", "Icelandic", "ISK (Icelandic króna)"),
    @("Switzerland", 8654622, 41284, 209.64, "Europe", "Western Europe", "French, Swiss German, Italian, Romansh", "CHF (Swiss franc)"),
    @("Christmas Island", 2072, 135, 15.35, "Oceania", "Australia and New Zealand", "English", "AUD (Australian dollar)"),
    @("Ireland", 4994724, 70273, 71.08, "Europe", "Northern Europe", "English, Irish", "EUR (Euro)"),
    @("Greenland", 56367, 2166086, 0.03, "Americas", "North America", "Greenlandic", "DKK (krone)"),
    @("Bouvet Island", 0, 49, 0, "Antarctic", "N/A", "Norwegian", "N/A"),
    @("New Zealand", 5084300, 270467, 18.8, "Oceania", "Australia and New Zealand", "English, Māori, New Zealand Sign Language", "NZD (New Zealand dollar)"),
    @("Norfolk Island", 2302, 36, 63.94, "Oceania", "Australia and New Zealand", "English, Norfuk", "AUD (Australian dollar)"),
    @("Thailand", 69799978, 513120, 136.03, "Asia", "South-Eastern Asia", "Thai", "THB (Thai baht)")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $ws.Cells.Item($row, 8).Value = $entry[7]
    $row++
}

# Remove the now-extra rows 13-16 (data previously present, now outside A1:H12 range)
$ws.Range("A13:H16").Delete()
